$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9 (Leve Item ID 5487)
$ws.Range("H9").Value = 504.83334
$ws.Range("I9").Value = 533.4545000000001
$ws.Range("K9").Value = 533.4545000000001
$ws.Range("M9").Value = -364.4545000000001
# Row 17 (Leve Item ID 38956)
$ws.Range("H17").Value = 2976
$ws.Range("J17").Value = 2976
$ws.Range("L17").Value = 8928
$ws.Range("N17").Value = -9264
# Row 19 (Leve Item ID 7015)
$ws.Range("H19").Value = 816.35297
$ws.Range("I19").Value = 920.1818
$ws.Range("J19").Value = 626
$ws.Range("K19").Value = 920.1818
$ws.Range("L19").Value = 626
$ws.Range("M19").Value = -745.1818
$ws.Range("N19").Value = -976
# Row 32 (Leve Item ID 5484)
$ws.Range("H32").Value = 3273.1667
$ws.Range("I32").Value = 2303.75
$ws.Range("K32").Value = 2303.75
$ws.Range("M32").Value = -1977.75
# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 3246.2334
$ws.Range("J40").Value = 5490.8184
$ws.Range("L40").Value = 5490.8184
$ws.Range("N40").Value = -5840.8184
# Row 69 (Leve Item ID 12616)
$ws.Range("H69").Value = 14899.857
$ws.Range("I69").Value = 14824.75
$ws.Range("K69").Value = 44474.25
$ws.Range("M69").Value = -43600.25
# Row 72 (Leve Item ID 12616)
$ws.Range("H72").Value = 14899.857
$ws.Range("I72").Value = 14824.75
$ws.Range("K72").Value = 133422.75
$ws.Range("M72").Value = -129054.75
# Row 87 (Leve Item ID 10651)
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
# Row 90 (Leve Item ID 10651)
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
# Row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 2284.2144
$ws.Range("I98").Value = 2118.3
$ws.Range("K98").Value = 2118.3
$ws.Range("M98").Value = -620.3000000000002
# Row 100 (Leve Item ID 19906)
$ws.Range("H100").Value = 2536.5
$ws.Range("I100").Value = 2536.5
$ws.Range("K100").Value = 2536.5
$ws.Range("M100").Value = -1995.5
# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 6213.4287
$ws.Range("I116").Value = 7249
$ws.Range("J116").Value = 4832.6665
$ws.Range("K116").Value = 7249
$ws.Range("L116").Value = 4832.6665
$ws.Range("M116").Value = -3807
$ws.Range("N116").Value = -11716.6665
# Row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 2284.2144
$ws.Range("I122").Value = 2118.3
$ws.Range("K122").Value = 6354.900000000001
$ws.Range("M122").Value = -3904.900000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 97 (Leve Item ID 19941)
$ws.Range("H97").Value = 2717.6
$ws.Range("J97").Value = 5447.5
$ws.Range("L97").Value = 5447.5
$ws.Range("N97").Value = -6439.5
# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 3999
$ws.Range("I102").Value = 3999
$ws.Range("K102").Value = 3999
$ws.Range("M102").Value = -2377
# Row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 8491.322
$ws.Range("I110").Value = 7188.5
$ws.Range("K110").Value = 7188.5
$ws.Range("M110").Value = -5143.5
# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 3986.75
$ws.Range("I122").Value = 4171.909
$ws.Range("J122").Value = 1950
$ws.Range("K122").Value = 12515.727
$ws.Range("L122").Value = 5850
$ws.Range("M122").Value = -10065.727
$ws.Range("N122").Value = -10750
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 4640.028
$ws.Range("I132").Value = 4057.5386
$ws.Range("K132").Value = 12172.6158
$ws.Range("M132").Value = -9642.6158

$ws = $wb.Worksheets.Item("BSM")
# Row 22 (Leve Item ID 5092)
$ws.Range("H22").Value = 600
$ws.Range("I22").Value = 600
$ws.Range("K22").Value = 600
$ws.Range("M22").Value = -427
# Row 76 (Leve Item ID 10630)
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
# Row 79 (Leve Item ID 10630)
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 2267.25
$ws.Range("I94").Value = 1675.2307
$ws.Range("J94").Value = 4832.6665
$ws.Range("K94").Value = 1675.2307
$ws.Range("L94").Value = 4832.6665
$ws.Range("M94").Value = -1224.2307
$ws.Range("N94").Value = -5734.6665
# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 1782.875
$ws.Range("I99").Value = 1077.3334
$ws.Range("K99").Value = 1077.3334
$ws.Range("M99").Value = 420.6666

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (Leve Item ID 5367)
$ws.Range("H22").Value = 40000000
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
# Row 86 (Leve Item ID 12584)
$ws.Range("H86").Value = 10619.8
$ws.Range("I86").Value = 4524.75
$ws.Range("K86").Value = 4524.75
$ws.Range("M86").Value = -3401.75
# Row 89 (Leve Item ID 12584)
$ws.Range("H89").Value = 10619.8
$ws.Range("I89").Value = 4524.75
$ws.Range("K89").Value = 22623.75
$ws.Range("M89").Value = -17007.75
# Row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 1540.2727
$ws.Range("I122").Value = 556.3333
$ws.Range("K122").Value = 1668.9999
$ws.Range("M122").Value = 781.0001
# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 1664.3077
$ws.Range("I134").Value = 1432.4
$ws.Range("J134").Value = 2437.3333
$ws.Range("K134").Value = 4297.200000000001
$ws.Range("L134").Value = 7311.999899999999
$ws.Range("M134").Value = -1762.200000000001
$ws.Range("N134").Value = -12381.9999

$ws = $wb.Worksheets.Item("CUL")
# Row 109 (Leve Item ID 27854)
$ws.Range("H109").Value = 3593.3333
$ws.Range("I109").Value = 2366.6667
$ws.Range("K109").Value = 7100.000100000001
$ws.Range("M109").Value = -6060.000100000001
# Row 110 (Leve Item ID 27857)
$ws.Range("H110").Value = 4400
$ws.Range("I110").Value = 4400
$ws.Range("K110").Value = 13200
$ws.Range("M110").Value = -9110
# Row 113 (Leve Item ID 27843)
$ws.Range("H113").Value = 2121.7144
$ws.Range("I113").Value = 2979.8333
$ws.Range("K113").Value = 8939.499899999999
$ws.Range("M113").Value = -6769.499899999999
# Row 138 (Leve Item ID 44105)
$ws.Range("H138").Value = 7573.091
$ws.Range("I138").Value = 7573.091
$ws.Range("K138").Value = 22719.273
$ws.Range("M138").Value = -17579.273

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 12616.167
$ws.Range("I70").Value = 11966
$ws.Range("K70").Value = 11966
$ws.Range("M70").Value = -11696
# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 12616.167
$ws.Range("I73").Value = 11966
$ws.Range("K73").Value = 11966
$ws.Range("M73").Value = -11030
# Row 97 (Leve Item ID 19940)
$ws.Range("H97").Value = 2630.0715
$ws.Range("I97").Value = 1934.625
$ws.Range("J97").Value = 3557.3333
$ws.Range("K97").Value = 1934.625
$ws.Range("L97").Value = 3557.3333
$ws.Range("M97").Value = -1438.625
$ws.Range("N97").Value = -4549.3333

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 2543.375
$ws.Range("I7").Value = 2169.4
$ws.Range("K7").Value = 2169.4
$ws.Range("M7").Value = -2057.4
# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 4180.4287
$ws.Range("I40").Value = 3852.8
$ws.Range("K40").Value = 3852.8
$ws.Range("M40").Value = -3716.8
# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 1436.421
$ws.Range("I46").Value = 1096.5
$ws.Range("K46").Value = 1096.5
$ws.Range("M46").Value = -908.5
# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 1147.5883
$ws.Range("I55").Value = 322.8
$ws.Range("K55").Value = 322.8
$ws.Range("M55").Value = -149.8
# Row 68 (Leve Item ID 12563)
$ws.Range("H68").Value = 7649.5
$ws.Range("I68").Value = 7649.5
$ws.Range("K68").Value = 7649.5
$ws.Range("M68").Value = -6900.5
# Row 71 (Leve Item ID 12563)
$ws.Range("H71").Value = 7649.5
$ws.Range("I71").Value = 7649.5
$ws.Range("K71").Value = 38247.5
$ws.Range("M71").Value = -34503.5
# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 3094.6
$ws.Range("I93").Value = 2915
$ws.Range("K93").Value = 2915
$ws.Range("M93").Value = -1667
# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 2543.375
$ws.Range("I126").Value = 2169.4
$ws.Range("K126").Value = 6508.200000000001
$ws.Range("M126").Value = -4038.200000000001
# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 18238.75
$ws.Range("I136").Value = 1526.75
$ws.Range("J136").Value = 68374.75
$ws.Range("K136").Value = 4580.25
$ws.Range("L136").Value = 205124.25
$ws.Range("M136").Value = -2030.25
$ws.Range("N136").Value = -210224.25

$ws = $wb.Worksheets.Item("WVR")
# Row 14 (Leve Item ID 2658)
$ws.Range("H14").Value = 5900
$ws.Range("I14").Value = 1000
$ws.Range("J14").Value = 9166.666999999999
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 9166.666999999999
$ws.Range("M14").Value = -832
$ws.Range("N14").Value = -9502.666999999999
# Row 96 (Leve Item ID 19977)
$ws.Range("H96").Value = 6523.6
$ws.Range("I96").Value = 6323
$ws.Range("K96").Value = 6323
$ws.Range("M96").Value = -4950
# Row 103 (Leve Item ID 18548)
$ws.Range("H103").Value = 40149.5
$ws.Range("J103").Value = 40149.5
$ws.Range("L103").Value = 40149.5
$ws.Range("N103").Value = -42493.5
# Row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 432.83334
$ws.Range("I113").Value = 337.5
$ws.Range("J113").Value = 528.1667
$ws.Range("K113").Value = 1012.5
$ws.Range("L113").Value = 1584.5001
$ws.Range("M113").Value = 1157.5
$ws.Range("N113").Value = -5924.5001
# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 11923.583
$ws.Range("I136").Value = 10887.5625
$ws.Range("K136").Value = 32662.6875
$ws.Range("M136").Value = -30112.6875
